# Add a new "14-10-2020" snapshot column (AC) to the active-cases sheet,
# mirroring the existing daily columns (the sheet already runs B..AB, one
# column per reporting date).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell AC1 --------------------------------------------------
# Clone AB1's formatting (bold, centered, thin border) onto AC1, then set
# the label as literal text so Excel doesn't reinterpret the dd-mm-yyyy
# string as a date serial (matches how the other header cells are stored).
$ws.Range("AB1").Copy($ws.Range("AC1"))
$ws.Range("AC1").Value = "14-10-2020"

# --- Data cells AC2:AC36 -----------------------------------------------
# Active-case counts per state/UT for the new date, same row order as the
# rest of the table (column AC = column index 29).
$values = @{
    2  = 199
    3  = 42855
    4  = 2960
    5  = 28897
    6  = 10835
    7  = 1127
    8  = 27208
    9  = 85
    10 = 21490
    11 = 4316
    12 = 15187
    13 = 10319
    14 = 2507
    15 = 9866
    16 = 7617
    17 = 113478
    18 = 95493
    19 = 969
    20 = 14661
    21 = 205884
    22 = 2867
    23 = 2367
    24 = 119
    25 = 1513
    26 = 22892
    27 = 4572
    28 = 8212
    29 = 21924
    30 = 344
    31 = 43239
    32 = 23728
    33 = 3500
    34 = 6576
    35 = 38082
    36 = 30988
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 29).Value = $values[$row]
}
